$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 13 (C13:J13) with the word "PLAYABLE" as an example board for testing purposes.
# E13 already contains "A"; the rest of the row was empty.
$ws.Range("C13").Value = "P"
$ws.Range("D13").Value = "L"
$ws.Range("E13").Value = "A"
$ws.Range("F13").Value = "Y"
$ws.Range("G13").Value = "A"
$ws.Range("H13").Value = "B"
$ws.Range("I13").Value = "L"
$ws.Range("J13").Value = "E"
